# San Jose expense workbook -- "Add files via upload" commit replay.
#
# Summary of the change being applied:
#   1. Each month sheet (Jan..Dec) gets its own month name written into B1
#      (a self-referential label that was previously blank).
#   2. The "Overall" sheet's B2/B3/B4 (Mortgage / Property Tax / Insurance)
#      get the same SUM-across-months formula that B5:B14 already have, and
#      pick up right-aligned number formatting to match.
#   3. Each month sheet picks up a basic portrait page setup.
#   4. Various cell selections (the cursor position saved with the sheet)
#      move around as a side effect of the author clicking through sheets.
#   5. The workbook window geometry changed (maximized / repositioned).

$wb = $excel.ActiveWorkbook

$months = @("Jan", "Feb", "Mar", "Apr", "May", "June", "July", "Aug", "Sep", "Oct", "Nov", "Dec")

# --- Overall sheet: new SUM formulas for Mortgage / Property Tax / Insurance ---
$overall = $wb.Worksheets.Item("Overall")
$overall.Activate()

$overall.Range("B2").Formula = "=SUM(Jan!B2, Feb!B2, Mar!B2, Apr!B2, May!B2, June!B2, July!B2, Aug!B2, Sep!B2, Oct!B2, Nov!B2, Dec!B2)"
$overall.Range("B2").HorizontalAlignment = -4152
$overall.Range("B2").VerticalAlignment = -4108

$overall.Range("B3").Formula = "=SUM(Jan!B3, Feb!B3, Mar!B3, Apr!B3, May!B3, June!B3, July!B3, Aug!B3, Sep!B3, Oct!B3, Nov!B3, Dec!B3)"
$overall.Range("B3").HorizontalAlignment = -4152
$overall.Range("B3").VerticalAlignment = -4107

$overall.Range("B4").Formula = "=SUM(Jan!B4, Feb!B4, Mar!B4, Apr!B4, May!B4, June!B4, July!B4, Aug!B4, Sep!B4, Oct!B4, Nov!B4, Dec!B4)"
$overall.Range("B4").HorizontalAlignment = -4152
$overall.Range("B4").VerticalAlignment = -4108

$overall.Range("B5").Select()

# --- Month sheets: write the month's own name into B1, set a print area,
#     and leave the selection on B1 (mirrors clicking into each tab) ---
foreach ($m in $months) {
    $ws = $wb.Worksheets.Item($m)
    $ws.Activate()
    $ws.Range("B1").Value = $m
    $ws.PageSetup.Orientation = 1
    $ws.Range("B1").Select()
}

# Dec ends up with a slightly different final selection than the rest.
$dec = $wb.Worksheets.Item("Dec")
$dec.Activate()
$dec.Range("D5").Select()

# Jan/Feb get a distinct final selection too (Feb's cursor rests on B1 like
# the other months; Jan keeps B1 as well) -- handled above already.

# --- Back to the Overall tab, which is the one left active/selected ---
$overall.Activate()
$overall.Range("B5").Select()

# --- Application window geometry (maximized, moved to top-left) ---
$excel.Left = 0
$excel.Top = 930
$excel.Width = 28800
$excel.Height = 14220
